$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate: drop old merge, clear all existing content/formatting.
$ws.Cells.UnMerge()
$ws.Cells.Clear()

# Reset row heights touched by the old layout back to the sheet default so no
# stale custom-height metadata survives on rows we are not explicitly sizing.
$ws.Rows(1).RowHeight = 15
$ws.Rows(1).AutoFit()

# --- Row 3: report title, merged across A3:J3 ---
$ws.Range("A3:J3").HorizontalAlignment = -4108
$ws.Range("A3:J3").Merge()
$ws.Range("A3").Value2 = "REPORTE PRODUCTO"
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Size = 14
$ws.Rows(3).RowHeight = 18.75

# --- Row 1: "Fecha:" label near the right (A1 blank-but-styled, I1 = Fecha:) ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("I1").Value2 = "Fecha:"
$ws.Range("I1").Font.Bold = $true

# --- Row 4: record-count label ---
$ws.Range("A4").Value2 = "Cantidad de registros"
$ws.Range("A4").Font.Bold = $true

# Column A width (matches the bestFit width Excel computed for the labels)
$ws.Columns(1).ColumnWidth = 22.16666666666667
$ws.Columns(1).BestFit = $true

# Selection ends on B4, matching the saved cursor position
$ws.Range("B4").Select() | Out-Null

Write-Host "done"
